$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 874.4
$ws.Range("I92").Value = 416
$ws.Range("K92").Value = 416
$ws.Range("M92").Value = 832
$ws.Range("H98").Value = 2815.2903
$ws.Range("I98").Value = 1529.7693
$ws.Range("J98").Value = 9500
$ws.Range("K98").Value = 1529.7693
$ws.Range("L98").Value = 9500
$ws.Range("M98").Value = -31.76929999999993
$ws.Range("N98").Value = -12496
$ws.Range("H122").Value = 2815.2903
$ws.Range("I122").Value = 1529.7693
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 4589.3079
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -2139.3079
$ws.Range("N122").Value = -33400
$ws.Range("H132").Value = 119904.68
$ws.Range("I132").Value = 126617.47
$ws.Range("K132").Value = 379852.41
$ws.Range("M132").Value = -377322.41
$ws.Range("H137").Value = 3500.8936
$ws.Range("I137").Value = 2930.5789
$ws.Range("J137").Value = 5908.8887
$ws.Range("K137").Value = 8791.736699999999
$ws.Range("L137").Value = 17726.6661
$ws.Range("M137").Value = -6241.736699999999
$ws.Range("N137").Value = -22826.6661
$ws.Range("H138").Value = 2726.05
$ws.Range("I138").Value = 1298.2307
$ws.Range("J138").Value = 2939.4023
$ws.Range("K138").Value = 3894.6921
$ws.Range("L138").Value = 8818.206900000001
$ws.Range("M138").Value = 1245.3079
$ws.Range("N138").Value = -19098.2069

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 27996
$ws.Range("J27").Value = 27996
$ws.Range("L27").Value = 27996
$ws.Range("N27").Value = -28364
$ws.Range("H63").Value = 5543540.5
$ws.Range("I63").Value = 11544701
$ws.Range("J63").Value = 4007.6924
$ws.Range("K63").Value = 11544701
$ws.Range("L63").Value = 4007.6924
$ws.Range("M63").Value = -11544015
$ws.Range("N63").Value = -5379.6924
$ws.Range("H66").Value = 5543540.5
$ws.Range("I66").Value = 11544701
$ws.Range("J66").Value = 4007.6924
$ws.Range("K66").Value = 57723505
$ws.Range("L66").Value = 20038.462
$ws.Range("M66").Value = -57720073
$ws.Range("N66").Value = -26902.462
$ws.Range("H132").Value = 2371.2856
$ws.Range("I132").Value = 1674.2
$ws.Range("J132").Value = 4114
$ws.Range("K132").Value = 5022.6
$ws.Range("L132").Value = 12342
$ws.Range("M132").Value = -2492.6
$ws.Range("N132").Value = -17402
$ws.Range("H137").Value = 41164.8
$ws.Range("J137").Value = 41164.8
$ws.Range("L137").Value = 41164.8
$ws.Range("N137").Value = -51364.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H99").Value = 2390.3333
$ws.Range("I99").Value = 1547.7778
$ws.Range("J99").Value = 3232.889
$ws.Range("K99").Value = 1547.7778
$ws.Range("L99").Value = 3232.889
$ws.Range("M99").Value = -49.77780000000007
$ws.Range("N99").Value = -6228.889
$ws.Range("H132").Value = 49648.57
$ws.Range("J132").Value = 49648.57
$ws.Range("L132").Value = 49648.57
$ws.Range("N132").Value = -59768.57
$ws.Range("H135").Value = 39275
$ws.Range("J135").Value = 39275
$ws.Range("L135").Value = 39275
$ws.Range("N135").Value = -49415
$ws.Range("H137").Value = 40614
$ws.Range("J137").Value = 40614
$ws.Range("L137").Value = 40614
$ws.Range("N137").Value = -50814
$ws.Range("H138").Value = 41240.57
$ws.Range("J138").Value = 41240.57
$ws.Range("L138").Value = 41240.57
$ws.Range("N138").Value = -51520.57
$ws.Range("H140").Value = 47419.8
$ws.Range("J140").Value = 47419.8
$ws.Range("L140").Value = 47419.8
$ws.Range("N140").Value = -57779.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 34940
$ws.Range("J23").Value = 34940
$ws.Range("L23").Value = 34940
$ws.Range("N23").Value = -35420
$ws.Range("H27").Value = 34940
$ws.Range("J27").Value = 34940
$ws.Range("L27").Value = 34940
$ws.Range("N27").Value = -35324
$ws.Range("H132").Value = 3562.5356
$ws.Range("I132").Value = 1650.8182
$ws.Range("K132").Value = 4952.4546
$ws.Range("M132").Value = -2422.4546
$ws.Range("H140").Value = 41147
$ws.Range("J140").Value = 41147
$ws.Range("L140").Value = 41147
$ws.Range("N140").Value = -51507

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2621.5344
$ws.Range("I121").Value = 214.14285
$ws.Range("J121").Value = 2951.9607
$ws.Range("K121").Value = 642.4285500000001
$ws.Range("L121").Value = 8855.882100000001
$ws.Range("M121").Value = 667.5714499999999
$ws.Range("N121").Value = -11475.8821
$ws.Range("H131").Value = 1017.6491
$ws.Range("J131").Value = 901.1177
$ws.Range("L131").Value = 2703.3531
$ws.Range("N131").Value = -12783.3531

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 40000
$ws.Range("J27").Value = 40000
$ws.Range("L27").Value = 40000
$ws.Range("N27").Value = -40332
$ws.Range("H135").Value = 43527390
$ws.Range("J135").Value = 43527390
$ws.Range("L135").Value = 43527390
$ws.Range("N135").Value = -43537530
$ws.Range("H137").Value = 74149.75
$ws.Range("J137").Value = 74149.75
$ws.Range("L137").Value = 74149.75
$ws.Range("N137").Value = -84349.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 29999
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 29999
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 29999
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -30589
$ws.Range("H93").Value = 1656.2632
$ws.Range("I93").Value = 569.5454999999999
$ws.Range("J93").Value = 3150.5
$ws.Range("K93").Value = 569.5454999999999
$ws.Range("L93").Value = 3150.5
$ws.Range("M93").Value = 678.4545000000001
$ws.Range("N93").Value = -5646.5
$ws.Range("H132").Value = 3785.7
$ws.Range("I132").Value = 1651.2646
$ws.Range("J132").Value = 5801.5557
$ws.Range("K132").Value = 4953.793799999999
$ws.Range("L132").Value = 17404.6671
$ws.Range("M132").Value = -2423.793799999999
$ws.Range("N132").Value = -22464.6671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3024.75
$ws.Range("I122").Value = 1883.64
$ws.Range("J122").Value = 5618.1816
$ws.Range("K122").Value = 5650.92
$ws.Range("L122").Value = 16854.5448
$ws.Range("M122").Value = -3200.92
$ws.Range("N122").Value = -21754.5448
$ws.Range("H132").Value = 11496507
$ws.Range("I132").Value = 818
$ws.Range("J132").Value = 19611110
$ws.Range("K132").Value = 2454
$ws.Range("L132").Value = 58833330
$ws.Range("M132").Value = 76
$ws.Range("N132").Value = -58838390
$ws.Range("H136").Value = 1474.44
$ws.Range("I136").Value = 740.73334
$ws.Range("J136").Value = 2575
$ws.Range("K136").Value = 2222.20002
$ws.Range("L136").Value = 7725
$ws.Range("M136").Value = 327.7999799999998
$ws.Range("N136").Value = -12825

Write-Host "Applied all Chocobo_Profits updates"